$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values on row 4 as per the 2024-10-10 FlashScore refresh.
$ws.Range("G4").Value = 2.82
$ws.Range("H4").Value = 3.05
$ws.Range("J4").Value = 3.3
$ws.Range("L4").Value = 2.95
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 1.91
$ws.Range("R4").Value = 1.8
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 14.5
$ws.Range("AA4").Value = 24
$ws.Range("AE4").Value = 12.5
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 13
$ws.Range("AL4").Value = 19.5
$ws.Range("AM4").Value = 27
$ws.Range("AN4").Value = 4.85
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 20
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 6.3
$ws.Range("AW4").Value = 4.45
$ws.Range("AY4").Value = 18
$ws.Range("BA4").Value = 70
